# Reorders/reshuffles the per-day subject-label cells inside the week's
# merged day-blocks (row groups A2:A7, A9:A14, A16:A21, A30:A35) and tweaks
# a handful of row heights to match the re-flowed layout.
#
# All the label cells share one cell style (yellow fill, centered /
# wrap-text) -- style index "4" in the original sheet. Rather than trying
# to reconstruct that formatting from scratch (fill color, alignment,
# wrap) we borrow it from a cell that keeps the same content + style in
# both the "before" and "after" layouts (F12 -- "Англ.яз. (Ольга)"), via
# Copy/PasteSpecial(xlPasteFormats). This guarantees the resulting style
# index is reused rather than a near-duplicate style being appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$donor = $ws.Range("F12")

function Set-LabelCell {
    param([string]$addr, [string]$text)
    $donor.Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).Value = $text
}

$excel.CutCopyMode = $false

# --- cells that lose their label entirely -------------------------------
$clearAddrs = @("C2","D2","F2","G2","G3","H3","F4","F11","D17","F17","F18")
foreach ($addr in $clearAddrs) {
    $ws.Range($addr).Clear()
}

# --- cells that gain / change a label -----------------------------------
Set-LabelCell "E2"  "Физика (ЮН)"
Set-LabelCell "C4"  "Химия (G1234re213)"
Set-LabelCell "D4"  "Англ.яз. (Ольга)"
Set-LabelCell "E4"  "Химия (Greeeg)"
Set-LabelCell "H4"  "Физика (ЮН)"
Set-LabelCell "F5"  "Англ.яз. (Ольга)"
Set-LabelCell "H5"  "Физика (ЮН)"
Set-LabelCell "C6"  "Физика (ЮН)"
Set-LabelCell "E6"  "Химия (Greeeg)"
Set-LabelCell "G6"  "География (Гриц)"
Set-LabelCell "H7"  "Физика (Гриц)"
Set-LabelCell "D10" "Англ.яз. (Ольга)"
Set-LabelCell "E10" "Химия (Greeeg)"
Set-LabelCell "H10" "Физика (Гриц)"
Set-LabelCell "C11" "Физика (ЮН)"
Set-LabelCell "E12" "Химия (Greeeg)"
Set-LabelCell "G13" "География (Гриц)"
Set-LabelCell "H13" "Физика (ЮН)"
Set-LabelCell "F14" "Англ.яз. (Ольга)"
Set-LabelCell "E17" "Химия (Greeeg)"
Set-LabelCell "H18" "Физика (Гриц)"
Set-LabelCell "G31" "География (Гриц)"

$excel.CutCopyMode = $false

# --- row-height touch-ups -------------------------------------------------
$ws.Rows.Item(2).RowHeight = 18
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 27
$ws.Rows.Item(7).RowHeight = 22
$ws.Rows.Item(10).RowHeight = 27
$ws.Rows.Item(11).RowHeight = 18
$ws.Rows.Item(13).RowHeight = 27
$ws.Rows.Item(14).RowHeight = 27
$ws.Rows.Item(17).RowHeight = 23
$ws.Rows.Item(18).RowHeight = 22
$ws.Rows.Item(31).RowHeight = 27
